# Add a "Relative" (column C) summary to each of the 7 datacompy report
# sheets: for rows 2-6 and 8 the relative value is B/B5 ("Total Rows"),
# for row 7 it's B7/B6 ("Source1 Rows") and for row 9 it's B9/B8
# ("Source2 Rows"). Division by zero is written as the literal text
# "<NaN>" (matching the desired fixture).  Also widen the new column,
# extend each sheet's AutoFilter from A1:B9 to A1:C9, and repoint the
# hidden _FilterDatabase defined names at the new range.

$wb = $excel.ActiveWorkbook

$sheetOrder = @("full","left","right","inner","complement","lcomp","rcomp")

$data = @{
    "full"       = @{ B = @(140,960,859,999,967,7,992,32); C = @(0.14014014014014015,0.960960960960961,0.8598598598598599,1.0,0.9679679679679679,0.007238883143743537,0.992992992992993,0.03225806451612903) }
    "left"       = @{ B = @(140,960,827,967,967,7,960,0);  C = @(0.14477766287487073,0.9927611168562565,0.8552223371251293,1.0,1.0,0.007238883143743537,0.9927611168562565,0.0) }
    "right"      = @{ B = @(140,960,852,992,960,0,992,32); C = @(0.14112903225806453,0.967741935483871,0.8588709677419355,1.0,0.967741935483871,0.0,1.0,0.03225806451612903) }
    "inner"      = @{ B = @(140,960,820,960,960,0,960,0);  C = @(0.14583333333333334,1.0,0.8541666666666666,1.0,1.0,0.0,1.0,0.0) }
    "complement" = @{ B = @(0,0,39,39,7,7,32,32);          C = @(0.0,0.0,1.0,1.0,0.1794871794871795,1.0,0.8205128205128205,1.0) }
    "lcomp"      = @{ B = @(0,0,7,7,7,7,0,0);              C = @(0.0,0.0,1.0,1.0,1.0,1.0,0.0,"<NaN>") }
    "rcomp"      = @{ B = @(0,0,32,32,0,0,32,32);          C = @(0.0,0.0,1.0,1.0,0.0,"<NaN>",1.0,1.0) }
}

foreach ($name in $sheetOrder) {
    $ws = $wb.Worksheets.Item($name)
    $rowVals = $data[$name]

    # --- column B corrections (only "inner" actually changes) -------------
    for ($i = 0; $i -lt 8; $i++) {
        $row = 2 + $i
        $ws.Range("B$row").Value = $rowVals.B[$i]
    }

    # --- header C1: "Relative", formatted like the existing B1 header -----
    $ws.Range("B1").Copy()
    $ws.Range("C1").PasteSpecial(-4122)
    $ws.Range("C1").Value = "Relative"

    # --- data cells C2:C9 ---------------------------------------------------
    for ($i = 0; $i -lt 8; $i++) {
        $row = 2 + $i
        $ws.Range("C$row").Value = $rowVals.C[$i]
    }
    $ws.Range("C2:C9").Interior.ColorIndex = 64
    $ws.Range("C2:C9").NumberFormat = "0.00%"

    # --- widen the new column ----------------------------------------------
    $ws.Columns.Item(3).ColumnWidth = 11.6

    # --- move the AutoFilter from A1:B9 to A1:C9 ----------------------------
    $ws.Range("A1:B9").AutoFilter()
    $ws.Range("A1:C9").AutoFilter()
}

# --- repoint the hidden _FilterDatabase defined names at the new range -----
foreach ($name in $sheetOrder) {
    $defName = $name + "!_FilterDatabase"
    $wb.Names.Item($defName).RefersTo = "=" + $name + "!`$A`$1:`$C`$9"
}

Write-Host "done"
